$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Row 8 (A8="Date", B8=generation timestamp)
$meta.Range("B8").Value = "2024-12-18T18:27:33+00:00"

# Row 22 (A22="Count", B22=number of concepts). The Count cell stores its
# value as text (like every other Property/Value pair on this sheet), so
# build the new text "8" on a scratch cell and paste just the value into
# B22 - a plain `.Value = "8"` assignment would be auto-converted to a
# number by Excel, which would change the stored cell type.
$scratch = $meta.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "8"
$scratch.Copy()
$meta.Range("B22").PasteSpecial(-4163)
$scratch.Clear()

# --- Concepts sheet: add a new concept row ---
$concepts = $wb.Worksheets.Item("Concepts")

# Duplicate the last data row (row 8) into the new row 9 so the new row
# inherits the same formatting/style without minting a new style entry,
# and column A ("Level") keeps its existing "1" shared-string value.
$concepts.Range("A8:D8").Copy($concepts.Range("A9:D9"))

# Overwrite the Code/Display columns for the new concept.
$concepts.Range("B9").Value = "Annotated-SNV"
$concepts.Range("C9").Value = "Annotated SNV"
